$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped crypto values.
# D-column price values are forced to text via NumberFormat "@" so Excel does not
# reinterpret numeric-looking strings (e.g. "307.69") as actual numbers, matching
# the workbook convention where these cells store text like "27.209.58".

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '27.209.58'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.67%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.901.11'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +0.63%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.9998'
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '307.69'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.48%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.5206'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +0.40%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3770'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +0.42%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.07276'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +1.02%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '21.17'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +0.17%  '

$ws.Range("E11").Value = '  +0.19%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.08289'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +8.66%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.908.17'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.97%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '96.39'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +2.00%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '5.276'
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.9998'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -0.13%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.000008633'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +1.42%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '14.57'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.81%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.9996'
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '27.232.17'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.56%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '5.091'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +0.78%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '2.157.64'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +1.95%  '

$ws.Range("E23").Value = '  +0.52%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '6.424'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.321'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +1.08%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '147.06'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.86%  '

$ws.Range("E27").Value = '  +0.41%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '18.20'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +0.69%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '115.07'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +0.79%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '4.837'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +0.94%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '4.903'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -0.20%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.09248'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +0.53%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.05072'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.64%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.7983'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +4.01%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.239'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -0.27%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '3.430'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +4.70%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.942'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -1.12%  '

$ws.Range("E38").Value = '  -0.08%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.5708'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +1.89%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.02003'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +0.78%  '

$ws.Range("E41").Value = '  +0.46%  '

$ws.Range("E42").Value = '  +0.09%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '6.582'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.66%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '116.69'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -1.95%  '

$ws.Range("E45").Value = '  +0.67%  '

$ws.Range("E46").Value = '  +0.59%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '10.09'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -0.46%  '

$ws.Range("E49").Value = '  +1.97%  '

$ws.Range("E50").Value = '  -0.15%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '63.89'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.09%  '
